$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 ---
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11

# --- Row 9 ---
$ws.Range("C9").Value = 16

# --- Row 10 ---
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

# --- Row 11 ---
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9

# --- Row 12 ---
$ws.Range("C12").Value = 10

# --- Row 13 ---
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

# --- Row 14 ---
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11

# --- Row 15 ---
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11

# --- New row 16 ---
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

# --- New row 17 ---
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true

# Copy the bold/bordered/centered style from A15 onto the two new A-column cells
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
